# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# For D-column values that look like plain decimals, Excel's COM layer would
# normally auto-coerce the assigned string into a floating point Double
# (dropping trailing zeros / introducing binary rounding noise such as
# 530.88999999999999). To keep those cells as literal text - matching how
# the sheet already stores every Price/Volume cell as a string - we flip the
# cell to Text format ("@") before writing the value, then restore the
# "Normal" style afterwards so no stray per-cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.772.48'
$ws.Range("E2").Value = '  +1.84%  '

$ws.Range("D3").Value = '3.166.84'
$ws.Range("E3").Value = '  +1.37%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.19%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '530.89'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '140.35'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.52%  '

$ws.Range("E7").Value = '  -0.04%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +10.90%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '7.28'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("E10").Value = '  +6.27%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.112'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +4.25%  '

$ws.Range("E12").Value = '  +2.41%  '

$ws.Range("D13").Value = '3.696.05'
$ws.Range("E13").Value = '  +1.03%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '25.78'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.71%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0000170'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +3.93%  '

$ws.Range("D16").Value = '58.683.89'
$ws.Range("E16").Value = '  +1.52%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '6.24'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.78%  '

$ws.Range("D18").Value = '3.150.21'
$ws.Range("E18").Value = '  +0.88%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.97'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.14%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '8.14'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.80%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '374.53'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.72%  '

$ws.Range("E22").Value = '  +0.38%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.531'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +4.99%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '69.65'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.77%  '

$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("E26").Value = '  +0.03%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.36'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +14.57%  '

$ws.Range("D28").Value = '0.0₃0859'
$ws.Range("E28").Value = '  -0.49%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '22.44'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +5.04%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.67%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '5.99'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.08%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.16'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.36%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.15'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.63%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '6.33'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +4.18%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '156.59'
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.34'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +3.34%  '

$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '25.25'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.79%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '2.708.92'
$ws.Range("E38").Value = '  +6.57%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.69'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.09%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0693'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.26%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.29'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +5.55%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0293'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +8.88%  '

$ws.Range("E43").Value = '  +2.98%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '39.11'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.47%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.36%  '

$ws.Range("D46").Value = '3.199.26'
$ws.Range("E46").Value = '  +1.05%  '

$ws.Range("E47").Value = '  +12.04%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '6.21'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.70%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.977'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.09%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '20.10'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.11%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.750'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.45%  '
